$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the post entry that used to occupy row 642 ("花より団子おおお（違う）...").
# Deleting the entire row shifts every row below it up by one, matching the
# renumbering seen across the rest of the sheet (old row 643 -> new row 642, etc.)
$ws.Rows.Item(642).Delete()
